$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

function Set-ParagraphText($index, $new) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    # Exclude the trailing paragraph mark so formatting / paragraph
    # structure is preserved; only the run text is replaced.
    $body = $d.Range($r.Start, $r.End - 1)
    $body.Text = $new
}

# Title: "Industrial Processes I" -> "Industrial Chemical Processes"
Replace-Text "Industrial Processes I" "Industrial Chemical Processes"

# Ativação date bump
Replace-Text "Ativação: 01/01/2024" "Ativação: 01/01/2025"

# Objetivos (PT)
Replace-Text "Conferir aos alunos uma visão geral da indústria química e correlatas, processos e produtos, e uma visão global das matérias primas mais importantes da indústria química." "Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos desta indústria."

# Objectives (EN, italic run)
Replace-Text "Objectives:Check the students an overview of the chemical industry and related industries, processes and products, and an overview of the most important raw materials in the chemical industry.." "Providing to the students an overview of the chemical and related industries, as well as the main features of the processes and production arrangements of this industry."

# Programa (PT) - long string (>255 chars), use paragraph Range instead of Find
Set-ParagraphText 14 "O conteúdo desta disciplina será de acordo com os tópicos a serem programados, devendo abordar assuntos relevantes relacionados a processos químicos e correlatas."

# Programa (EN, italic) - long string (>255 chars), use paragraph Range instead of Find
Set-ParagraphText 15 "The content of this subject will be in accordance with the topics to be programmed, and must address relevant subjects related to chemical and related processes."

# Método
Replace-Text "Aulas expositivas, filmes e leituras de artigos técnicos" "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

# Critério
Replace-Text "Provas em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula e frequência." "A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:NF = NOTA x % FREQ."

# Norma de recuperação
Replace-Text "Prova escrita para alunos que tenham média final maior ou igual a 3,0 (Três) e inferior a 5,0 (Cinco). A nota final será a média aritmética entre a média final e a prova escrita." "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos."

# Bibliografia - long string (>255 chars), use paragraph Range instead of Find
Set-ParagraphText 19 "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual Econômico da Indústria Química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;"

Write-Output "Done"
